$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "warden 1" row (row 2), shifting all subsequent rows up.
$ws.Rows.Item(2).Delete()
